$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" (E1) and "is_enabled_lbl" (F1) template columns
# and shift the remaining columns (dept_ids_lbl, role_ids_lbl, rem) left.
$ws.Range("E1:F1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
